# Churn project results_final.xlsx update
# Re-ranks a few classifiers (column A) and refreshes the metrics
# (columns E:L) for every classifier row on the results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (rank) changes ---------------------------------------------
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 12
$ws.Range("A6").Value = 7

# --- Column B (classifier label) changes ----------------------------------
$ws.Range("B3").Value = "Random Forest"
$ws.Range("B4").Value = "XGB"
$ws.Range("B6").Value = "Stacking (SGD)"

# --- Columns E:L (metrics) refresh for every data row ----------------------
$metrics = @{
    2  = @(0.9742554265522464, 0.9634517766497462, 0.9844398340248963, 0.9937435019264876, 949, 981, 15, 36)
    3  = @(0.9581019687026754, 0.9461928934010152, 0.9688149688149689, 0.9925932155016003, 932, 966, 30, 53)
    4  = @(0.9727410398788491, 0.9614213197969543, 0.9833852544132918, 0.9923888447189774, 947, 980, 16, 38)
    5  = @(0.9540636042402827, 0.9451776649746193, 0.9617768595041323, 0.9912370293356166, 931, 959, 37, 54)
    6  = @(0.9742554265522464, 0.9644670050761421, 0.9834368530020704, 0.9879997145944183, 950, 980, 16, 35)
    7  = @(0.9727410398788491, 0.9604060913705583, 0.9843912591050988, 0.9760279697470082, 946, 981, 15, 39)
    8  = @(0.9570923775870772, 0.950253807106599,  0.9629629629629629, 0.9570546143966729, 936, 960, 36, 49)
    9  = @(0.8642099949520444, 0.8812182741116751, 0.8509803921568627, 0.9167273153527817, 868, 844, 152, 117)
    10 = @(0.9121655729429581, 0.9218274111675127, 0.9034825870646767, 0.9122189264672904, 908, 899, 97, 77)
    11 = @(0.9040888440181727, 0.9045685279187817, 0.9027355623100304, 0.9043870915132611, 891, 900, 96, 94)
    12 = @(0.8182735991923271, 0.7756345177664975, 0.8460686600221484, 0.8563553707214646, 764, 857, 139, 221)
    13 = @(0.7895002523977789, 0.8040609137055837, 0.7795275590551181, 0.8407406274845576, 792, 772, 224, 193)
    14 = @(0.6456335184250379, 0.867005076142132,  0.5992982456140351, 0.7009061627219538, 854, 425, 571, 131)
}

foreach ($row in $metrics.Keys) {
    $vals = $metrics[$row]
    $ws.Range("E$row").Value = $vals[0]
    $ws.Range("F$row").Value = $vals[1]
    $ws.Range("G$row").Value = $vals[2]
    $ws.Range("H$row").Value = $vals[3]
    $ws.Range("I$row").Value = $vals[4]
    $ws.Range("J$row").Value = $vals[5]
    $ws.Range("K$row").Value = $vals[6]
    $ws.Range("L$row").Value = $vals[7]
}
